$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Time Period" row (row 8): end month moves from 2024-08 to 2024-09 ---
$ws.Range("B8:F8").Value = "2003-05:2024-09"

# --- Update the "Update" row (row 10): refresh date moves from 2024-09-20 to 2024-10-19 ---
# Leading apostrophe forces text (matches the source file's quote-prefixed text cells,
# otherwise Excel would reinterpret the string as a date serial number).
$ws.Range("B10:F10").Value = "'2024-10-19"

# --- Append the new monthly data point (row 261) for 2024-09-30 ---
# Copy number formats from the last existing data row so the new row keeps the same
# date / number styles (yyyy-mm date style, #,##0.00 style) instead of General.
$ws.Range("A260:F260").Copy()
$ws.Range("A261:F261").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A261").Value = 45565
$ws.Range("B261").Value = 74094
$ws.Range("C261").Value = 1035
$ws.Range("D261").Value = 47385
$ws.Range("E261").Value = 13953
$ws.Range("F261").Value = 11721

# --- Row 1 re-measures to a 15pt height once the workbook is resaved (the red Calibri
# font used by the Wind-link cell has a taller natural line height than the sheet's
# default font) ---
$ws.Rows.Item(1).RowHeight = 15
